# Apply "table output of vacancies" update to the "Статистика по годам" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header labels (engineer -> driver)
$ws.Range("C1").Value = "Средняя зарплата - водитель"
$ws.Range("E1").Value = "Количество вакансий - водитель"

# Widen columns C and E by 1 unit
# Note: the ColumnWidth COM property adds a fixed padding offset (5/6)
# relative to the raw OOXML column "width" attribute when saved, so the
# input value is adjusted here to land exactly on width=29 / width=32.
$ws.Columns.Item(3).ColumnWidth = 28.166666666666668
$ws.Columns.Item(5).ColumnWidth = 31.166666666666668

# Update the data values for the new profession ("водитель")
$values = @(
    @(2, 56419, 123),
    @(3, 65786, 696),
    @(4, 64078, 792),
    @(5, 66416, 1162),
    @(6, 73525, 1434),
    @(7, 76036, 1557),
    @(8, 78612, 1918),
    @(9, 80402, 1935),
    @(10, 91682, 2059),
    @(11, 78512, 2604),
    @(12, 93515, 2534),
    @(13, 95512, 3050),
    @(14, 104304, 3225),
    @(15, 111293, 3123),
    @(16, 119943, 1723),
    @(17, 142799, 525)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
}
